$wb = $excel.ActiveWorkbook

# 1) Rename "Full Frame" -> "Full-Frame" on the main table sheet.
$ws1 = $wb.Worksheets.Item("BigAssTableOfImagingFormats")
$ws1.Range("A2").Value = "Full-Frame"

# 2) Add a new Changelog entry describing the rename.
$ws2 = $wb.Worksheets.Item("Changelog")
$ws2.Range("A4").Value = 42815
$ws2.Range("A4").NumberFormat = $ws2.Range("A3").NumberFormat
$ws2.Range("B4").Value = 'Changed "Full Frame" to "Full-Frame" to match BATOL'
